# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a 🚀
#
# Updates the "Metadata" sheet (Version/Date/Publisher/Jurisdiction rows,
# removing the stray duplicated "Contact" row) and the "Elements" sheet
# (root Extension row's Short/Definition text) to match the republished
# StructureDefinition-practitioner-hierarchy-level-description content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Metadata": Property/Value table
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bumped to the new publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$meta.Range("B9").Value = "Alvearie Team"

# The old row 10 "Contact" / "No display for ContactDetail" becomes
# "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# remove it entirely so everything below shifts up by one.
$meta.Rows.Item(11).Delete()

# ---------------------------------------------------------------------
# Sheet "Elements": StructureDefinition element table
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; its Short (K) / Definition (L)
# text is updated from the generic placeholder to the real description.
$elements.Range("K2").Value = "Practitioner Hierarchy Level Description"
$elements.Range("L2").Value = "Description of the level of the practitioner within the organinzational hierarchy"

# The "Short" column (K) widened slightly to fit the new text (the sheet
# uses best-fit/auto-size columns, so its stored width shifts with the
# new longest value).
$elements.Columns.Item(11).ColumnWidth = 36.8
